$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Wnt9a from ECs): updated ligand expression values from new TPM data,
# plus the downstream derived-specificity / edge-weight columns that are
# recalculated from them.
$ws.Range("G2").Value = 0.1887043333333333
$ws.Range("H2").Value = 0.566113
$ws.Range("I2").Value = 0.02109097403787168
$ws.Range("J2").Value = 0.02109097403787168
$ws.Range("Q2").Value = 0.009049630812222223
$ws.Range("R2").Value = 0.08144667730999999
$ws.Range("S2").Value = 0.02109097403787168
$ws.Range("T2").Value = 0.02109097403787168

# Row 3 (Wnt9a from FAPs): recalculated derived-specificity values plus a
# tiny floating point refresh of G3/R3.
$ws.Range("G3").Value = 7.8617
$ws.Range("I3").Value = 0.8786809908633213
$ws.Range("J3").Value = 0.8786809908633211
$ws.Range("R3").Value = 3.393188337
$ws.Range("S3").Value = 0.8786809908633213
$ws.Range("T3").Value = 0.8786809908633211

# Row 4 (Wnt9a from MuSCs): recalculated derived-specificity values.
$ws.Range("I4").Value = 0.1002280350988072
$ws.Range("J4").Value = 0.1002280350988071
$ws.Range("S4").Value = 0.1002280350988072
$ws.Range("T4").Value = 0.1002280350988071
